$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 29, pushing the existing rows 29-46 down to 31-48.
$ws.Rows("29:30").Insert()

# Row 29: new weekly entry (Especial), date 2023-03-24 (serial 45009)
$ws.Range("A29").Value = 8
$ws.Range("B29").Value = "Terminal La Palmera de La Serena"
$ws.Range("C29").Value = "Coquimbo"
$ws.Range("D29").Value = 45009
$ws.Range("E29").Value = 4
$ws.Range("F29").Value = "Fruta"
$ws.Range("G29").Value = 100107
$ws.Range("H29").Value = "Otros"
$ws.Range("I29").Value = 100107011
$ws.Range("J29").Value = "Tuna"
$ws.Range("K29").Value = "Sin especificar"
$ws.Range("L29").Value = "Especial"
$ws.Range("M29").Value = 300
$ws.Range("N29").Value = 13000
$ws.Range("O29").Value = 14000
$ws.Range("P29").Value = 13500
$ws.Range("Q29").Value = "`$/caja 18 kilos"
$ws.Range("R29").Value = "Provincia de Limarí"
$ws.Range("S29").Value = 750
$ws.Range("T29").Value = 18

# Row 30: new weekly entry (Primera), same date
$ws.Range("A30").Value = 8
$ws.Range("B30").Value = "Terminal La Palmera de La Serena"
$ws.Range("C30").Value = "Coquimbo"
$ws.Range("D30").Value = 45009
$ws.Range("E30").Value = 4
$ws.Range("F30").Value = "Fruta"
$ws.Range("G30").Value = 100107
$ws.Range("H30").Value = "Otros"
$ws.Range("I30").Value = 100107011
$ws.Range("J30").Value = "Tuna"
$ws.Range("K30").Value = "Sin especificar"
$ws.Range("L30").Value = "Primera"
$ws.Range("M30").Value = 200
$ws.Range("N30").Value = 11000
$ws.Range("O30").Value = 12000
$ws.Range("P30").Value = 11500
$ws.Range("Q30").Value = "`$/caja 18 kilos"
$ws.Range("R30").Value = "Provincia de Limarí"
$ws.Range("S30").Value = 639
$ws.Range("T30").Value = 18
